$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 594; everything currently at/after
# row 594 shifts down by one (row 652 -> 653).
$ws.Rows.Item(594).Insert()

# Populate the newly inserted row 594 with the new weekly data point.
$ws.Cells.Item(594, 1).Value = 4
$ws.Cells.Item(594, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(594, 3).Value = "Los Lagos"
$ws.Cells.Item(594, 4).Value = 45132
$ws.Cells.Item(594, 5).Value = 10
$ws.Cells.Item(594, 6).Value = "Fruta"
$ws.Cells.Item(594, 7).Value = 100102
$ws.Cells.Item(594, 8).Value = "Cítricos"
$ws.Cells.Item(594, 9).Value = 100102006
$ws.Cells.Item(594, 10).Value = "Pomelo"
$ws.Cells.Item(594, 11).Value = "Start Ruby"
$ws.Cells.Item(594, 12).Value = "Primera"
$ws.Cells.Item(594, 13).Value = 100
$ws.Cells.Item(594, 14).Value = 14000
$ws.Cells.Item(594, 15).Value = 14000
$ws.Cells.Item(594, 16).Value = 14000
$ws.Cells.Item(594, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(594, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(594, 19).Value = 1000
$ws.Cells.Item(594, 20).Value = 14
